$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new centered, bold title paragraph at the very start of the
#    document: "MOTOR FAILURE 02 (BALANCED) - SPARK"
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$titleRange = $firstPara.Range
$titleRange.Collapse(1)  # wdCollapseStart
$titleRange.InsertParagraphBefore()

# The newly created (now first) paragraph is empty; fill it in.
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Text = [char]0x2013
$titlePara.Range.Text = "MOTOR FAILURE 02 (BALANCED) " + [char]0x2013 + " SPARK"
$titlePara.Alignment = 1  # wdAlignParagraphCenter
$titlePara.Range.Font.Bold = 1
$titlePara.Range.Font.BoldBi = 1
$titlePara.Range.Font.Size = 12
$titlePara.Range.Font.SizeBi = 12

Write-Host "Paragraph count after insert:" $d.Paragraphs.Count
Write-Host "Title text: [$($d.Paragraphs.Item(1).Range.Text)]"
